# Compas sheet: insert a new header/title row above the existing header,
# add a merged "Formations" / "Diplomes" banner, adjust the related
# structures (autofilter, named range, conditional formatting, comments)
# that reference row coordinates, and a couple of small cosmetic tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compas")
$ws.Activate()

# ---------------------------------------------------------------------
# 1) Insert a new row 1. Everything that was on row 1/2 shifts to 2/3.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).Insert()

# ---------------------------------------------------------------------
# 2) Build the new banner row (row 1).
# ---------------------------------------------------------------------

# A1:H1 -> big merged, shaded title band (currently left blank)
$bandA = $ws.Range("A1:H1")
$bandA.Font.Bold = $true
$bandA.Font.Name = "Arial"
$bandA.HorizontalAlignment = -4108   # xlCenter
$bandA.VerticalAlignment = -4108     # xlCenter
$bandA.Interior.Color = 15064278     # D6DCE5 (theme dk2 @ 80% tint)
$bandA.Borders.Item(8).Weight = -4138  # xlEdgeTop, xlMedium
$bandA.Borders.Item(9).Weight = -4138  # xlEdgeBottom, xlMedium
$ws.Range("A1").Borders.Item(7).Weight = -4138  # xlEdgeLeft
$ws.Range("H1").Borders.Item(10).Weight = -4138 # xlEdgeRight
$bandA.Merge()

# J1:K1 -> "Formations" (reuses the blue shade used lower for FC1/FC2)
$bandJ = $ws.Range("J1:K1")
$bandJ.NumberFormat = "m/d/yyyy"
$bandJ.Interior.Color = 12176622    # BDD7EE (theme accent5 @ 60% tint)
$bandJ.Font.Bold = $true
$bandJ.Font.Name = "Arial"
$bandJ.HorizontalAlignment = -4108
$bandJ.VerticalAlignment = -4108
$bandJ.WrapText = $true
$bandJ.Borders.Item(8).Weight = -4138
$bandJ.Borders.Item(9).Weight = -4138
$ws.Range("J1").Borders.Item(7).Weight = -4138
$ws.Range("K1").Borders.Item(10).Weight = -4138
$bandJ.Merge()
$ws.Range("J1").Value = "Formations"

# M1:N1 -> "Diplomes" (reuses the green shade used lower for AFPS/PSC1)
$bandM = $ws.Range("M1:N1")
$bandM.NumberFormat = "m/d/yyyy"
$bandM.Interior.Color = 11848645    # C5E0B4 (theme accent6 @ 60% tint)
$bandM.Font.Bold = $true
$bandM.Font.Name = "Arial"
$bandM.HorizontalAlignment = -4108
$bandM.VerticalAlignment = -4108
$bandM.WrapText = $true
$bandM.Borders.Item(8).Weight = -4138
$bandM.Borders.Item(9).Weight = -4138
$ws.Range("M1").Borders.Item(7).Weight = -4138
$ws.Range("N1").Borders.Item(10).Weight = -4138
$bandM.Merge()
$ws.Range("M1").Value = "Diplomes"

$ws.Rows.Item(1).RowHeight = 15.75

Write-Host "Banner row built"

# ---------------------------------------------------------------------
# 3) Column A is very slightly narrower now.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.43

# ---------------------------------------------------------------------
# 4) The filter/database range now starts on row 2 (the real header).
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Compas!_FilterDatabase") {
        $n.RefersTo = "=Compas!`$A`$2:`$R`$3"
    }
}

$ws.AutoFilterMode = $false
$ws.Range("A2:R3").AutoFilter()

# ---------------------------------------------------------------------
# 5) Conditional formatting ranges shift down by one row too.
# ---------------------------------------------------------------------
$fcsQR = $ws.Range("Q2:R65535").FormatConditions
for ($i = 1; $i -le $fcsQR.Count; $i++) {
    $fcsQR.Item($i).ModifyAppliesToRange($ws.Range("Q3:R65536"))
}

$fcsG = $ws.Range("G2:G65535").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fc = $fcsG.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("G3:G65536"))
    $fc.Formula1 = $fc.Formula1.Replace("G2", "G3")
}

$fcsH = $ws.Range("H2:H65535").FormatConditions
for ($i = 1; $i -le $fcsH.Count; $i++) {
    $fc = $fcsH.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("H3:H65536"))
    $fc.Formula1 = $fc.Formula1.Replace("H2", "H3")
}

Write-Host "Ranges updated"

# ---------------------------------------------------------------------
# 6) Comments were anchored to row 1; move them down to row 2 along
#    with the header cells they annotate.
# ---------------------------------------------------------------------
$commentCells = @("J", "K", "M", "N", "P", "Q", "R")
foreach ($col in $commentCells) {
    $oldCell = $ws.Range($col + "1")
    $newCell = $ws.Range($col + "2")
    $cmt = $oldCell.Comment
    if ($cmt -ne $null) {
        $txt = $cmt.Text()
        $cmt.Delete()
        $newCell.AddComment($txt)
    }
}

Write-Host "Comments relocated"
